$wb = $excel.ActiveWorkbook

# --- Update timestamps on the "data" sheet (column F, rows 2-8) ---
$dataSheet = $wb.Worksheets.Item("data")
$dataSheet.Range("F2").Value = "2021-10-05 14:33:51.789283"
$dataSheet.Range("F3").Value = "2021-10-05 14:33:51.789292"
$dataSheet.Range("F4").Value = "2021-10-05 14:33:51.789296"
$dataSheet.Range("F5").Value = "2021-10-05 14:33:51.789299"
$dataSheet.Range("F6").Value = "2021-10-05 14:33:51.789302"
$dataSheet.Range("F7").Value = "2021-10-05 14:33:51.789305"
$dataSheet.Range("F8").Value = "2021-10-05 14:33:51.789308"

# --- Add the new "metadata" sheet positioned after "data" ---
$meta = $wb.Worksheets.Add($null, $dataSheet)
$meta.Name = "metadata"

# Header row (row 1, columns B..G)
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row (row 2)
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Frontonasal dysplasia"
$meta.Range("C2").Value = 104
$meta.Range("E2").Value = "2021-01-18T23:49:34.519092Z"
$meta.Range("F2").Value = "2021-10-05 14:33:51.785459"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/104/?format=json"

# D2 must stay literal text "1.0" (not get auto-coerced to the number 1).
# Build it as a text-returning formula on a scratch cell, then paste-special
# the computed value back in as a static value so it is no longer a formula.
$meta.Range("Z1").Formula = "=""1.0"""
$meta.Range("Z1").Copy()
$meta.Range("D2").PasteSpecial(-4163)
$meta.Range("Z1").Clear()
$excel.CutCopyMode = $false

# --- Match header / index-column styling to the "data" sheet (bold, bordered) ---
# Use copy / paste-special-formats so we reuse the existing style record
# instead of minting a new one.
$dataSheet.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$dataSheet.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the originally active sheet ("data") — adding the new sheet
# shifts focus onto it by default.
$dataSheet.Activate() | Out-Null
